# Weekly fruit/vegetable price update: insert two new "Choclo" price rows
# (dated 44946, the newest week) at the top of the data table, pushing the
# existing rows down by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right before the current row 110 (the first data
# row after the header + the 108/109 pair already covering that date).
$ws.Rows("110:111").Insert()

# New row 110: Choclero / Primera
$ws.Range("A110").Value = 11
$ws.Range("B110").Value = "Vega Monumental Concepción"
$ws.Range("C110").Value = "Bíobío"
$ws.Range("D110").Value = 44946
$ws.Range("E110").Value = 8
$ws.Range("F110").Value = 100112024
$ws.Range("G110").Value = "Choclo"
$ws.Range("H110").Value = "Choclero"
$ws.Range("I110").Value = "Primera"
$ws.Range("J110").Value = 5000
$ws.Range("K110").Value = 200
$ws.Range("L110").Value = 200
$ws.Range("M110").Value = 200
$ws.Range("N110").Value = "$/unidad"
$ws.Range("O110").Value = "Región Metropolitana"
$ws.Range("P110").Value = 200
$ws.Range("Q110").Value = 1
$ws.Range("R110").Value = "Hortaliza"

# New row 111: Choclero / Segunda
$ws.Range("A111").Value = 11
$ws.Range("B111").Value = "Vega Monumental Concepción"
$ws.Range("C111").Value = "Bíobío"
$ws.Range("D111").Value = 44946
$ws.Range("E111").Value = 8
$ws.Range("F111").Value = 100112024
$ws.Range("G111").Value = "Choclo"
$ws.Range("H111").Value = "Choclero"
$ws.Range("I111").Value = "Segunda"
$ws.Range("J111").Value = 3000
$ws.Range("K111").Value = 150
$ws.Range("L111").Value = 150
$ws.Range("M111").Value = 150
$ws.Range("N111").Value = "$/unidad"
$ws.Range("O111").Value = "Región Metropolitana"
$ws.Range("P111").Value = 150
$ws.Range("Q111").Value = 1
$ws.Range("R111").Value = "Hortaliza"
